# Update rules and advisors for the cryptic eye box
# Adds 8 new advisor rows (rows 66-73) to the "advisors" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A name, B reputation_cost, C coin_cost, D guild, E buy_goods,
# F sell_goods, G tax, H raid, I explore, J research, K upgrade, L build,
# M set, N glory, O stage, P ability

$rows = @(
    @{ A = "The Bringer of Storms";   B = 0; C = 6; D = "explorers"; E = 0; F = 0; G = 0; H = 0; I = 3; J = 0; K = 0;   L = 0;  M = "cryptic eye"; N = 1; O = 6; P = "When Bringer of Storms becomes your active advisor, exhaust an advisor ina nother province's council room. You may explore any dangerous waters space to recover a lost relic. The defense of the endeavor is 8 and it is a dangerous site. You gain glory as normal for the endeavor." },
    @{ A = "The Gatherer of Brethren"; B = 0; C = 6; D = "soldiers";  E = 0; F = 0; G = 0; H = 4; I = 0; J = 0; K = 0;   L = 0;  M = "cryptic eye"; N = 1; O = 6; P = "When Gatherer of Brethren becomes your active advisor, shine the light of truth on advisors in the forum. You must hire all society members (including public ones) for free. If you successfully raid a ship, you may gain a ship upgrade worth 10 gold from the supply." },
    @{ A = "The Holder of Secrets";    B = 0; C = 6; D = "soldiers";  E = 0; F = 0; G = 0; H = 3; I = 0; J = 0; K = 0;   L = 0;  M = "cryptic eye"; N = 1; O = 6; P = "When Holder of Secrets becomes your active advisor, you must take an advisor from a province's council room (if it doesn't have enmity on it). Place two enmity tokens on that council room. If you raid a ship carrying a relic, +3 dice to the endeavor." },
    @{ A = "The Keeper of Lore";       B = 0; C = 6; D = "builders";  E = 0; F = 0; G = 0; H = 0; I = 0; J = 0; K = -10; L = 0;  M = "cryptic eye"; N = 1; O = 6; P = "When Keeper of Lore becomes your active advisor, draw the top four cards of the Research deck and keep two." },
    @{ A = "The Master of Whispers";   B = 0; C = 6; D = "merchants"; E = 3; F = 0; G = 0; H = 0; I = 0; J = 0; K = 0;   L = 0;  M = "cryptic eye"; N = 1; O = 6; P = "When Masters of Whispers becomes your active advisor, you must take one gold from each other province's vault. Any goods you buy this turn immediately go to any warehouse you control." },
    @{ A = "The Seeker of Answers";    B = 0; C = 6; D = "builders";  E = 0; F = 0; G = 0; H = 0; I = 0; J = 0; K = 0;   L = -2; M = "cryptic eye"; N = 1; O = 6; P = "When Seeker of Answers becomes your active advisor, you must take an advisor from a province's council room (if it doesn't have enmity on it). Place two enmity tokens on that council room. You may use as many goods as you want to discount a build action this turn." },
    @{ A = "The Ninth Initiate";       B = 0; C = 6; D = "merchants"; E = 0; F = 3; G = 0; H = 0; I = 0; J = 0; K = 0;   L = 0;  M = "cryptic eye"; N = 1; O = 6; P = "When the Ninth Initiate becomes your active advisor, shine the light of truth on advisors in the forum. You must hire all private society members for free. You may take a card from any treasure room that doesn't have enmity on it. Place four enmity on that treasure room. This is not a raid or endeavor." },
    @{ A = "The Walker of Paths";      B = 0; C = 6; D = "explorers"; E = 0; F = 0; G = 0; H = 0; I = 3; J = 0; K = 0;   L = 0;  M = "cryptic eye"; N = 1; O = 6; P = "When Walker of Paths becomes your active advisor, shine the light of truth on advisors in the forum. You may hire all society members (including public ones) for a total of one reputation. You may explore any atoll to recover a lost relic. The defense of the endeavor is 8 and it is a dangerous site. You gain glory for the endeavor." }
)

$startRow = 66
$r = $startRow
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value  = $row.A
    $ws.Cells.Item($r, 2).Value  = $row.B
    $ws.Cells.Item($r, 3).Value  = $row.C
    $ws.Cells.Item($r, 4).Value  = $row.D
    $ws.Cells.Item($r, 5).Value  = $row.E
    $ws.Cells.Item($r, 6).Value  = $row.F
    $ws.Cells.Item($r, 7).Value  = $row.G
    $ws.Cells.Item($r, 8).Value  = $row.H
    $ws.Cells.Item($r, 9).Value  = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $r++
}

# Update selection to match the new bottom of the data (after the last row)
$endRow = $r
$ws.Range("A$endRow").Select()
